$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Force the cell to stay a literal (shared-string) text value even when
    # the text looks numeric (e.g. "1", "179"). Plain `.Value = "179"` gets
    # auto-coerced to a number by Excel's normal type inference, so we
    # briefly mark the cell as Text, assign it, then clear the number
    # format again so the cell keeps using the default style (s="0"),
    # matching the original workbook's formatting.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2 (employee MNV=2)
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "dương"
Set-TextValue $ws.Range("C2") "1"
Set-TextValue $ws.Range("D2") "1"
$ws.Range("E2").Value = 1

# Row 3 (employee MNV=5)
$ws.Range("A3").Value = 5
$ws.Range("B3").Value = "huy"
Set-TextValue $ws.Range("C3") "179"
$ws.Range("D3").Value = "@gmail.com"
$ws.Range("E3").Value = 9

# Row 4 (employee MNV=6)
$ws.Range("A4").Value = 6
$ws.Range("B4").Value = "huy"
Set-TextValue $ws.Range("C4") "179"
$ws.Range("D4").Value = "@gmail.com"
$ws.Range("E4").Value = 9

# Row 5 (employee MNV=7)
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "afafaf"
$ws.Range("C5").Value = "255252sffs"
$ws.Range("D5").Value = "ssfsfsf"
$ws.Range("E5").Value = 837002627
